$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D ("Tipo" shifts from D to E)
$ws.Columns.Item(4).Insert()

# New header for column D - copy formatting from the "R2" header cell, then set the text
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)
$ws.Cells.Item(1, 4).Value = "MAE"

# MAE values for rows 2-14
$maeValues = @(
    0.9239883288152645,
    0.8429226117853587,
    0.886280647950258,
    0.8956978851546247,
    1.027519701295929,
    1.350102492671558,
    2.046624516121589,
    2.091210877959376,
    0.6472536726993832,
    0.3447233267553575,
    0.8583964064824784,
    1.863313287418611,
    1.29915171380136
)

for ($i = 0; $i -lt $maeValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $maeValues[$i]
}
